# Auto-generated edit script applying the Sargatanas_Profits.xlsx diff
# Updates currentAveragePrice/currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ (M) and LeveProfitHQ (N) columns across all 8 sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 14).ClearContents()
$ws.Cells.Item(34, 8).Value = 1829.2858
$ws.Cells.Item(34, 9).Value = 467.5
$ws.Cells.Item(34, 11).Value = 467.5
$ws.Cells.Item(34, 13).Value = -264.5
$ws.Cells.Item(36, 8).Value = 1829.2858
$ws.Cells.Item(36, 9).Value = 467.5
$ws.Cells.Item(36, 11).Value = 467.5
$ws.Cells.Item(36, 13).Value = 247.5
$ws.Cells.Item(43, 8).Value = 458183
$ws.Cells.Item(43, 9).Value = 524.5
$ws.Cells.Item(43, 11).Value = 524.5
$ws.Cells.Item(43, 13).Value = -455.5
$ws.Cells.Item(53, 8).Value = 7875.375
$ws.Cells.Item(53, 9).Value = 10292.333
$ws.Cells.Item(53, 10).Value = 6425.2
$ws.Cells.Item(53, 11).Value = 10292.333
$ws.Cells.Item(53, 12).Value = 6425.2
$ws.Cells.Item(53, 13).Value = -9655.333000000001
$ws.Cells.Item(53, 14).Value = -7699.2
$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 14).ClearContents()
$ws.Cells.Item(64, 8).Value = 28578070
$ws.Cells.Item(64, 9).Value = 37043460
$ws.Cells.Item(64, 10).Value = 7374.75
$ws.Cells.Item(64, 11).Value = 37043460
$ws.Cells.Item(64, 12).Value = 7374.75
$ws.Cells.Item(64, 13).Value = -37043212
$ws.Cells.Item(64, 14).Value = -7870.75
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 14).ClearContents()
$ws.Cells.Item(67, 8).Value = 28578070
$ws.Cells.Item(67, 9).Value = 37043460
$ws.Cells.Item(67, 10).Value = 7374.75
$ws.Cells.Item(67, 11).Value = 37043460
$ws.Cells.Item(67, 12).Value = 7374.75
$ws.Cells.Item(67, 13).Value = -37042602
$ws.Cells.Item(67, 14).Value = -9090.75
$ws.Cells.Item(113, 8).Value = 60011190
$ws.Cells.Item(113, 9).Value = 3666
$ws.Cells.Item(113, 10).Value = 83347460
$ws.Cells.Item(113, 11).Value = 3666
$ws.Cells.Item(113, 12).Value = 83347460
$ws.Cells.Item(113, 13).Value = -412
$ws.Cells.Item(113, 14).Value = -83353968
$ws.Cells.Item(116, 8).Value = 10425751
$ws.Cells.Item(116, 9).Value = 19236848
$ws.Cells.Item(116, 11).Value = 19236848
$ws.Cells.Item(116, 13).Value = -19233406
$ws.Cells.Item(137, 8).Value = 3242.2222
$ws.Cells.Item(137, 10).Value = 3988.6
$ws.Cells.Item(137, 12).Value = 11965.8
$ws.Cells.Item(137, 14).Value = -17065.8
$ws.Cells.Item(141, 8).Value = 3752.4783
$ws.Cells.Item(141, 9).Value = 2851.2778
$ws.Cells.Item(141, 11).Value = 8553.8334
$ws.Cells.Item(141, 13).Value = -3373.8334

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(7, 8).Value = 60001
$ws.Cells.Item(7, 10).Value = 60001
$ws.Cells.Item(7, 12).Value = 60001
$ws.Cells.Item(7, 14).Value = -60229
$ws.Cells.Item(32, 8).Value = 3928583.8
$ws.Cells.Item(32, 9).Value = 4260612.5
$ws.Cells.Item(32, 10).Value = 27249.5
$ws.Cells.Item(32, 11).Value = 4260612.5
$ws.Cells.Item(32, 12).Value = 27249.5
$ws.Cells.Item(32, 13).Value = -4260325.5
$ws.Cells.Item(32, 14).Value = -27823.5
$ws.Cells.Item(61, 8).Value = 34485644
$ws.Cells.Item(61, 9).Value = 1242.6666
$ws.Cells.Item(61, 11).Value = 1242.6666
$ws.Cells.Item(61, 13).Value = -1030.6666
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 14).ClearContents()
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 4021.8918
$ws.Cells.Item(122, 9).Value = 2912.4443
$ws.Cells.Item(122, 11).Value = 8737.332900000001
$ws.Cells.Item(122, 13).Value = -6287.332900000001
$ws.Cells.Item(132, 8).Value = 3918.7314
$ws.Cells.Item(132, 9).Value = 3067.2654
$ws.Cells.Item(132, 11).Value = 9201.796200000001
$ws.Cells.Item(132, 13).Value = -6671.796200000001
$ws.Cells.Item(136, 8).Value = 34485644
$ws.Cells.Item(136, 9).Value = 1242.6666
$ws.Cells.Item(136, 11).Value = 3727.9998
$ws.Cells.Item(136, 13).Value = -1177.9998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 45455060
$ws.Cells.Item(80, 9).Value = 83333950
$ws.Cells.Item(80, 10).Value = 393
$ws.Cells.Item(80, 11).Value = 83333950
$ws.Cells.Item(80, 12).Value = 393
$ws.Cells.Item(80, 13).Value = -83332952
$ws.Cells.Item(80, 14).Value = -2389
$ws.Cells.Item(83, 8).Value = 45455060
$ws.Cells.Item(83, 9).Value = 83333950
$ws.Cells.Item(83, 10).Value = 393
$ws.Cells.Item(83, 11).Value = 416669750
$ws.Cells.Item(83, 12).Value = 1965
$ws.Cells.Item(83, 13).Value = -416664758
$ws.Cells.Item(83, 14).Value = -11949
$ws.Cells.Item(107, 8).Value = 187530820
$ws.Cells.Item(107, 9).Value = 225027980
$ws.Cells.Item(107, 11).Value = 225027980
$ws.Cells.Item(107, 13).Value = -225026060

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).ClearContents()
$ws.Cells.Item(4, 14).ClearContents()
$ws.Cells.Item(86, 8).Value = 54215536
$ws.Cells.Item(86, 9).Value = 53059650
$ws.Cells.Item(86, 11).Value = 53059650
$ws.Cells.Item(86, 13).Value = -53058527
$ws.Cells.Item(89, 8).Value = 54215536
$ws.Cells.Item(89, 9).Value = 53059650
$ws.Cells.Item(89, 11).Value = 265298250
$ws.Cells.Item(89, 13).Value = -265292634
$ws.Cells.Item(107, 8).Value = 3131.2727
$ws.Cells.Item(107, 9).Value = 2977.1428
$ws.Cells.Item(107, 11).Value = 2977.1428
$ws.Cells.Item(107, 13).Value = -1057.1428
$ws.Cells.Item(110, 8).Value = 62000
$ws.Cells.Item(110, 10).Value = 62000
$ws.Cells.Item(110, 12).Value = 62000
$ws.Cells.Item(110, 14).Value = -70180

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(94, 8).Value = 5800
$ws.Cells.Item(94, 10).Value = 5800
$ws.Cells.Item(94, 12).Value = 17400
$ws.Cells.Item(94, 14).Value = -18752
$ws.Cells.Item(137, 8).Value = 92490.13
$ws.Cells.Item(137, 9).Value = 67516.53
$ws.Cells.Item(137, 10).Value = 146005
$ws.Cells.Item(137, 11).Value = 202549.59
$ws.Cells.Item(137, 12).Value = 438015
$ws.Cells.Item(137, 13).Value = -197449.59
$ws.Cells.Item(137, 14).Value = -448215

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 2000183.6
$ws.Cells.Item(2, 9).Value = 110.85714
$ws.Cells.Item(2, 10).Value = 6667020
$ws.Cells.Item(2, 11).Value = 110.85714
$ws.Cells.Item(2, 12).Value = 6667020
$ws.Cells.Item(2, 13).Value = 2.142859999999999
$ws.Cells.Item(2, 14).Value = -6667246
$ws.Cells.Item(70, 8).Value = 8169.7334
$ws.Cells.Item(70, 9).Value = 4782.6
$ws.Cells.Item(70, 11).Value = 4782.6
$ws.Cells.Item(70, 13).Value = -4512.6
$ws.Cells.Item(73, 8).Value = 8169.7334
$ws.Cells.Item(73, 9).Value = 4782.6
$ws.Cells.Item(73, 11).Value = 4782.6
$ws.Cells.Item(73, 13).Value = -3846.6
$ws.Cells.Item(107, 8).Value = 1601840
$ws.Cells.Item(107, 9).Value = 2667066.8
$ws.Cells.Item(107, 10).Value = 4000
$ws.Cells.Item(107, 11).Value = 2667066.8
$ws.Cells.Item(107, 12).Value = 4000
$ws.Cells.Item(107, 13).Value = -2665146.8
$ws.Cells.Item(107, 14).Value = -7840
$ws.Cells.Item(132, 8).Value = 1543.0625
$ws.Cells.Item(132, 9).Value = 1160.7858
$ws.Cells.Item(132, 10).Value = 4219
$ws.Cells.Item(132, 11).Value = 3482.3574
$ws.Cells.Item(132, 12).Value = 12657
$ws.Cells.Item(132, 13).Value = -952.3574000000003
$ws.Cells.Item(132, 14).Value = -17717
$ws.Cells.Item(133, 8).Value = 56749
$ws.Cells.Item(133, 10).Value = 56749
$ws.Cells.Item(133, 12).Value = 56749
$ws.Cells.Item(133, 14).Value = -66869

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 39737.5
$ws.Cells.Item(2, 10).Value = 30483.334
$ws.Cells.Item(2, 12).Value = 30483.334
$ws.Cells.Item(2, 14).Value = -30707.334
$ws.Cells.Item(46, 8).Value = 2656.6072
$ws.Cells.Item(46, 9).Value = 2161.0557
$ws.Cells.Item(46, 10).Value = 3548.6
$ws.Cells.Item(46, 11).Value = 2161.0557
$ws.Cells.Item(46, 12).Value = 3548.6
$ws.Cells.Item(46, 13).Value = -1973.0557
$ws.Cells.Item(46, 14).Value = -3924.6
$ws.Cells.Item(122, 8).Value = 4411.609
$ws.Cells.Item(122, 9).Value = 3722.1
$ws.Cells.Item(122, 11).Value = 11166.3
$ws.Cells.Item(122, 13).Value = -8716.299999999999
$ws.Cells.Item(132, 8).Value = 10425194
$ws.Cells.Item(132, 9).Value = 23812680
$ws.Cells.Item(132, 11).Value = 71438040
$ws.Cells.Item(132, 13).Value = -71435510

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 762.56525
$ws.Cells.Item(113, 9).Value = 764.1739
$ws.Cells.Item(113, 11).Value = 2292.5217
$ws.Cells.Item(113, 13).Value = -122.5217000000002
$ws.Cells.Item(122, 8).Value = 26532154
$ws.Cells.Item(122, 9).Value = 45821776
$ws.Cells.Item(122, 11).Value = 137465328
$ws.Cells.Item(122, 13).Value = -137462878
